# Custom Question Generation Code
# Adds a new "Code" question (C1 - Java) row to Sheet2 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# --- Column D width (new column used by the coding question) ---
$ws.Columns.Item(4).ColumnWidth = 18.67

# --- Row 6 values, added in the same order the original author typed them ---
# Identifying / scoring columns (A-C)
$ws.Cells.Item(6, 1).Value = "C1"
$ws.Cells.Item(6, 2).Value = "Java"
$ws.Cells.Item(6, 3).Value = 30

# Question + constraints (D-G)
$ws.Cells.Item(6, 4).Value = "Write a function that takes a list of integers and returns the second largest number in the list. If there is no second largest (e.g., all numbers are equal), return None."
$ws.Cells.Item(6, 5).Value = "The list must contain at least two integers"
$ws.Cells.Item(6, 6).Value = "The integers can be negative or positive"
$ws.Cells.Item(6, 7).Value = "Do not use built-in sort functions"

# Example 1 (H-I)
$ws.Cells.Item(6, 8).Value = "Input: [10, 20, 4, 45, 99] "
$ws.Cells.Item(6, 9).Value = "Output: 45"

# Example 2 (K-L) typed before going back to fill in the explanations (J, M)
$ws.Cells.Item(6, 11).Value = "Input: [5, 5, 5]  "
$ws.Cells.Item(6, 12).Value = "Output: null"

# Explanations (J for example 1, M for example 2)
$ws.Cells.Item(6, 10).Value = "The largest number is 99`nThe second largest is 45"
$ws.Cells.Item(6, 13).Value = "All numbers are the same, so there is no distinct second largest"

# Points column (N)
$ws.Cells.Item(6, 14).Value = 70

# --- Formatting for row 6 ---
# A:C and N -> vertical top aligned, no wrap
$topAligned = $ws.Range("A6:C6")
$topAligned.VerticalAlignment = -4160
$ws.Range("N6").VerticalAlignment = -4160

# D:M -> vertical top aligned, wrap text
$wrapped = $ws.Range("D6:M6")
$wrapped.VerticalAlignment = -4160
$wrapped.WrapText = $true

# Row height
$ws.Rows.Item(6).RowHeight = 135

# --- Selection ---
$ws.Range("M6").Select() | Out-Null
